# Cambio diapo y Avd por Av.
#
# Slide 2 ("Problemática"): shorten paragraphs 2 and 3 of the content
# placeholder.
# Slide 3 ("Objetivos"): fix "Avd." -> "Av." typo.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 2 - shape "Marcador de contenido 6"
# ---------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$shape2 = $slide2.Shapes.Item(2)
$tr2 = $shape2.TextFrame.TextRange

$oldPara2 = "Se puede identificar un problema en común en cuanto tránsito y transporte, esto puede ser producto del crecimiento del número de vehículos, la falta de infraestructura vial, el mal estado de la misma y un crecimiento urbano descontrolado. "
$newPara2 = "Se puede identificar un problema en común en cuanto tránsito y transporte."
$found2 = $tr2.Find($oldPara2)
if ($found2 -ne $null) {
    $found2.Text = $newPara2
}

$oldPara3 = "Las congestiones vehiculares producidas en las horas pico (horas de mayor flujo vehicular), se deben a la necesidad de las personas de llegar a sus sitios de trabajo, vivienda o estudio, estos embotellamientos se presentan debido al gran número de vehículos que circulan por las vías y la falta de planes de movilidad, ocasionando malestar entre los conductores, usuarios de transporte público y peatones."
$newPara3 = "Las congestiones vehiculares producidas en las horas pico (horas de mayor flujo vehicular), se deben a la necesidad de las personas de llegar a sus sitios de trabajo, vivienda o estudio."
$found3 = $tr2.Find($oldPara3)
if ($found3 -ne $null) {
    $found3.Text = $newPara3
}

# ---------------------------------------------------------------------
# Slide 3 - shape "Marcador de contenido 6" : "Avd." -> "Av."
# ---------------------------------------------------------------------
$slide3 = $p.Slides.Item(3)
$shape3 = $slide3.Shapes.Item(2)
$tr3 = $shape3.TextFrame.TextRange

$oldAvd = "Tomar datos reales de la situación actual sobre el entorno de la congestión vehicular en la calle Alonso de Mercadillo, la Avd. Universitaria y la Av. Manuel Agustín Aguirre."
$newAvd = "Tomar datos reales de la situación actual sobre el entorno de la congestión vehicular en la calle Alonso de Mercadillo, la Av. Universitaria y la Av. Manuel Agustín Aguirre."
$foundAvd = $tr3.Find($oldAvd)
if ($foundAvd -ne $null) {
    $foundAvd.Text = $newAvd
} else {
    # Fallback: replace just the "Avd." occurrence if the whole-paragraph
    # match above didn't hit (keeps the edit robust to minor drift).
    $foundAvd2 = $tr3.Find("Avd.")
    if ($foundAvd2 -ne $null) {
        $foundAvd2.Text = "Av."
    }
}
